# Insert a new data row before row 304 (shifts existing rows 304..409 down to 305..410)
# and populate the new row 304 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, pushing rows 304:409 down to 305:410
$ws.Rows("304:304").Insert()

# Populate the newly inserted row 304 with the new record
$ws.Cells.Item(304, 1).Value  = 10
$ws.Cells.Item(304, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(304, 3).Value  = "La Araucanía"
$ws.Cells.Item(304, 4).Value  = 44988
$ws.Cells.Item(304, 5).Value  = 9
$ws.Cells.Item(304, 6).Value  = 100112001
$ws.Cells.Item(304, 7).Value  = "Berenjena"
$ws.Cells.Item(304, 8).Value  = "Sin especificar"
$ws.Cells.Item(304, 9).Value  = "Primera"
$ws.Cells.Item(304, 10).Value = 80
$ws.Cells.Item(304, 11).Value = 12000
$ws.Cells.Item(304, 12).Value = 12000
$ws.Cells.Item(304, 13).Value = 12000
$ws.Cells.Item(304, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(304, 15).Value = "Región del Maule"
$ws.Cells.Item(304, 16).Value = 300
$ws.Cells.Item(304, 17).Value = 40
$ws.Cells.Item(304, 18).Value = "Hortaliza"
